# Updates cryptos list values (Price / Volume(1h)) to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.819.71"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "'1.644.37"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "'216.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("E6").Value = "  -0.50%  "

$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("D10").Value = "'19.20"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.35%  "

$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("D12").Value = "'1.643.01"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.47%  "

$ws.Range("E13").Value = "  -0.77%  "

$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").Value = "'26.824.07"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("E17").Value = "  -1.38%  "

$ws.Range("D18").Value = "'214.79"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("E20").Value = "  +1.05%  "

$ws.Range("E21").Value = "  +9.48%  "

$ws.Range("E22").Value = "  -0.44%  "

$ws.Range("E23").Value = "  -1.88%  "

$ws.Range("D24").Value = "'146.44"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("E25").Value = "  +0.19%  "

$ws.Range("D26").Value = "'0.118"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.34%  "

$ws.Range("D27").Value = "'7.19"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("D28").Value = "'15.67"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.72%  "

$ws.Range("D29").Value = "'0.0508"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.94%  "

$ws.Range("E30").Value = "  +0.31%  "

$ws.Range("E31").Value = "  -1.11%  "

$ws.Range("D32").Value = "'3.01"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.21%  "

$ws.Range("D33").Value = "'1.287.84"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.95%  "

$ws.Range("E35").Value = "  +1.35%  "

$ws.Range("E36").Value = "  -1.26%  "

$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("E38").Value = "  -1.06%  "

$ws.Range("E39").Value = "  +0.36%  "

$ws.Range("D40").Value = "'0.808"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.23%  "

$ws.Range("D41").Value = "'2.23"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("E42").Value = "  -2.42%  "

$ws.Range("D43").Value = "'1.784.69"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("D44").Value = "'61.38"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.57%  "

$ws.Range("D45").Value = "'92.05"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.24%  "

$ws.Range("E46").Value = "  +1.38%  "

$ws.Range("E47").Value = "  -1.24%  "

$ws.Range("D48").Value = "'0.0521"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.74%  "

$ws.Range("D49").Value = "'7.66"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.74%  "

$ws.Range("D50").Value = "'0.0969"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("E51").Value = "  -0.03%  "
